$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1713
$ws1.Range("F4").Value = 1155
$ws1.Range("F9").Value = 95
$ws1.Range("F14").Value = 465
$ws1.Range("F18").Value = 690
$ws1.Range("F24").Value = 285
$ws1.Range("F36").Value = 23

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 718
$ws2.Range("F5").Value = 609
$ws2.Range("F6").Value = 609
$ws2.Range("F12").Value = 271
$ws2.Range("F15").Value = 339
$ws2.Range("F16").Value = 339
$ws2.Range("F19").Value = 929
$ws2.Range("F24").Value = 25
$ws2.Range("F26").Value = 229

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2263
$ws3.Range("F9").Value = 1118
$ws3.Range("F10").Value = 242
$ws3.Range("F11").Value = 74

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2263
$ws4.Range("F5").Value = 1713
$ws4.Range("F9").Value = 1119
$ws4.Range("F10").Value = 242
$ws4.Range("F11").Value = 74
$ws4.Range("F12").Value = 718
$ws4.Range("F13").Value = 1155
$ws4.Range("F17").Value = 609
$ws4.Range("F19").Value = 95
$ws4.Range("F25").Value = 465
$ws4.Range("F28").Value = 690
$ws4.Range("F32").Value = 285
$ws4.Range("F33").Value = 271
$ws4.Range("F39").Value = 339
$ws4.Range("F45").Value = 25
$ws4.Range("F46").Value = 229
